# Bulk Upload Template-v2.xlsx — "Complete bulk upload for role and technology"
#
# The Employees sheet's sample/demo data is trimmed down from 3 rows to a
# single example row, and that row's sample data is refreshed with a new
# EMP ID and a new email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Employees")

# --- Update the remaining sample row (row 2): new EMP ID + new email ---
$ws.Range("B2").Value = "E0133"

$g2 = $ws.Range("G2")
$g2.Value = "shonamishra170@gmail.com"

# --- Hyperlinks: drop the ones tied to the rows being removed (G3, G4) and
#     refresh the one that stays (G2) to point at the new email address.
#     (Individual Hyperlink.Delete() isn't wired up in this host, so the
#     collection is rebuilt from scratch.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($g2, "mailto:shonamishra170@gmail.com", "", "", "shonamishra170@gmail.com")

# --- Remove the two extra sample employees entirely (rows 3 and 4) ---
$ws.Rows("3:4").Delete()

# --- The row-delete ripples all the way to the bottom of the sheet (as it
#     would in a full 1,048,576-row grid), so the last two rows pick up the
#     sheet's post-delete row height too.
$ws.Cells.Item(1048575, 10).NumberFormat = "General"
$ws.Cells.Item(1048576, 10).NumberFormat = "General"
$ws.Rows(1048575).RowHeight = 12.8
$ws.Rows(1048576).RowHeight = 12.8

# --- Restore the selection cursor to where it ended up after the edit ---
$ws.Range("C6").Select()
